# 20230112 Pyspark test dataset
# Insert a new record (vinod / kadam) as the second data row (new row 4),
# pushing the existing data down by one row, then fix up the hyperlink
# collection (which the engine does not auto-shift on row insert) and the
# active selection the way Excel leaves it after such an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new row, shifting rows 4:24 down to 5:25 -----------------
$ws.Rows("4:4").Insert()

# --- 2. Populate the newly inserted row with the new record ----------------
$ws.Range("A4").Value = "vinod"
$ws.Range("B4").Value = "kadam"
$ws.Range("C4").Value = 31
$ws.Range("D4").Value = "m"
$ws.Range("E4").Value = "vinod#gmail.com"
$ws.Range("F4").Value = "borivali"
$ws.Range("G4").Value = "teacher"
$ws.Range("H4").Value = 2.8
$ws.Range("E4").Style = "Hyperlink"

# --- 3. Rebuild the hyperlink collection in the new row order --------------
# (the engine does not shift existing Hyperlink objects when rows are
# inserted, so we recreate them against the now-correct ranges)
$ws.Hyperlinks.Delete()

$mailRefs = @("E2","E3","E5","E6","E7","E8","E9","E10","E11","E12","E13","E14","E15","E16","E17","E18","E19","E20","E21","E22","E23","E24","E25")
$mailAddrs = @(
    "mailto:kiran@gmail.com",
    "mailto:kapil@gmail.com",
    "mailto:samira@yahoo.com",
    "mailto:vidhan@hotmail.com",
    "mailto:abhijit@gmail.com",
    "mailto:manrata@yahoo.com",
    "mailto:neha@hotmail.com",
    "mailto:shubham@hotmail.com",
    "mailto:darshan@gmail.com",
    "mailto:anuj@yahoo.com",
    "mailto:harshali@gmail.com",
    "mailto:nitesh@gmail.com",
    "mailto:nitesh_a@gmail.com",
    "mailto:vidit@yahoo.com",
    "mailto:shraddha@yahoo.com",
    "mailto:jinal@gmail.com",
    "mailto:rohit@gmail.com",
    "mailto:siddhesh@yahoo.com",
    "mailto:nidhi@hotmail.cpm",
    "mailto:harsh@yahoo.com",
    "mailto:snehal@gmail.com",
    "mailto:medha@yahoo.com",
    "mailto:sonal@gmail.com"
)

for ($i = 0; $i -lt $mailRefs.Length; $i++) {
    $ws.Hyperlinks.Add($ws.Range($mailRefs[$i]), $mailAddrs[$i])
    $ws.Range($mailRefs[$i]).Style = "Hyperlink"
}

# --- 4. Match the selection Excel leaves behind after this edit ------------
$ws.Range("H4").Select()
